# Refresh the cryptocurrency price/volume snapshot (cap3k "cryptos" sheet).
# Column D = Price, Column E = Volume(1h) change, both stored as text.
#
# Values that look like plain decimals (e.g. "321.41") would otherwise be
# auto-converted to numbers by Excel's input parser, which would both change
# the stored cell type and silently drop formatting such as trailing zeros
# ("9.00" -> 9). A leading apostrophe forces Excel to keep the literal text,
# exactly preserving the source string - the same trick a human typing these
# values into a spreadsheet would use.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.276.08"
$ws.Range("E2").Value = "  +2.25%  "

$ws.Range("D3").Value = "2.509.65"
$ws.Range("E3").Value = "  +1.00%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'321.41"

$ws.Range("D6").Value = "'108.72"
$ws.Range("E6").Value = "  +0.66%  "

$ws.Range("E7").Value = "  +1.14%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  +1.14%  "

$ws.Range("E10").Value = "  +2.25%  "

$ws.Range("D11").Value = "'20.09"
$ws.Range("E11").Value = "  +9.02%  "

$ws.Range("E12").Value = "  +1.13%  "

$ws.Range("E13").Value = "  +0.40%  "

$ws.Range("E14").Value = "  +0.59%  "

$ws.Range("D15").Value = "2.907.24"
$ws.Range("E15").Value = "  +1.21%  "

$ws.Range("D16").Value = "2.516.69"
$ws.Range("E16").Value = "  +1.29%  "

$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("D18").Value = "48.116.65"
$ws.Range("E18").Value = "  +2.09%  "

$ws.Range("E19").Value = "  -2.74%  "

$ws.Range("D20").Value = "'6.73"
$ws.Range("E20").Value = "  +1.78%  "

$ws.Range("E21").Value = "  +0.86%  "

$ws.Range("E22").Value = "  -1.26%  "

$ws.Range("D23").Value = "'72.19"
$ws.Range("E23").Value = "  +2.42%  "

$ws.Range("D24").Value = "'277.81"
$ws.Range("E24").Value = "  +12.93%  "

$ws.Range("D25").Value = "'2.55"
$ws.Range("E25").Value = "  +0.62%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "'25.83"

$ws.Range("E28").Value = "  +5.18%  "

$ws.Range("D29").Value = "'9.83"
$ws.Range("E29").Value = "  -1.09%  "

$ws.Range("D30").Value = "'35.56"
$ws.Range("E30").Value = "  +2.89%  "

$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("D32").Value = "'49.17"
$ws.Range("E32").Value = "  -1.22%  "

$ws.Range("D33").Value = "'19.59"
$ws.Range("E33").Value = "  -3.60%  "

$ws.Range("E34").Value = "  +0.85%  "

$ws.Range("E36").Value = "  +0.55%  "

$ws.Range("E37").Value = "  +0.65%  "

$ws.Range("E38").Value = "  -2.94%  "

$ws.Range("E39").Value = "  +0.92%  "

$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("D41").Value = "'122.21"
$ws.Range("E41").Value = "  +1.04%  "

$ws.Range("E42").Value = "  +0.43%  "

$ws.Range("D43").Value = "'21.65"
$ws.Range("E43").Value = "  -5.41%  "

$ws.Range("E44").Value = "  +3.56%  "

$ws.Range("D45").Value = "2.002.04"

$ws.Range("D46").Value = "'3.18"
$ws.Range("E46").Value = "  +5.62%  "

$ws.Range("E47").Value = "  +3.33%  "

$ws.Range("E48").Value = "  -0.86%  "

$ws.Range("D49").Value = "'9.00"
$ws.Range("E49").Value = "  -1.03%  "

$ws.Range("E50").Value = "  +2.66%  "

$ws.Range("D51").Value = "'80.28"
$ws.Range("E51").Value = "  +3.79%  "
